$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.072.61'
$ws.Range("E2").Value = '  +1.35%  '
$ws.Range("D3").Value = '3.530.29'
$ws.Range("E3").Value = '  +0.34%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '601.70'
$ws.Range("E5").Value = '  +1.33%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '184.56'
$ws.Range("E6").Value = '  +6.05%  '
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.596'
$ws.Range("E8").Value = '  +0.07%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.140'
$ws.Range("E9").Value = '  +4.63%  '
$ws.Range("E10").Value = '  -1.78%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.447'
$ws.Range("E11").Value = '  +2.38%  '
$ws.Range("D12").Value = '4.143.74'
$ws.Range("E12").Value = '  +0.46%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '32.63'
$ws.Range("E13").Value = '  +12.22%  '
$ws.Range("E14").Value = '  -0.31%  '
$ws.Range("D15").Value = '68.039.53'
$ws.Range("E15").Value = '  +1.36%  '
$ws.Range("E16").Value = '  +0.45%  '
$ws.Range("D17").Value = '3.529.22'
$ws.Range("E17").Value = '  +1.47%  '
$ws.Range("E18").Value = '  +1.51%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '14.89'
$ws.Range("E19").Value = '  +4.34%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '400.15'
$ws.Range("E20").Value = '  +1.29%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '8.14'
$ws.Range("E21").Value = '  +1.61%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '73.97'
$ws.Range("E22").Value = '  +1.19%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.548'
$ws.Range("E23").Value = '  +1.23%  '
$ws.Range("B25").Value = 'PEPE'
$ws.Range("C25").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0000125'
$ws.Range("E25").Value = '  +2.35%  '
$ws.Range("B26").Value = 'LEO'
$ws.Range("C26").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '5.69'
$ws.Range("E26").Value = '  -0.08%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.67'
$ws.Range("E27").Value = '  +3.39%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.179'
$ws.Range("E28").Value = '  -0.99%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.997'
$ws.Range("E29").Value = '  -0.12%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.33'
$ws.Range("E30").Value = '  +0.73%  '
$ws.Range("E31").Value = '  +1.80%  '
$ws.Range("E32").Value = '  +1.00%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '24.21'
$ws.Range("E33").Value = '  +1.21%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '7.46'
$ws.Range("E34").Value = '  +1.18%  '
$ws.Range("E35").Value = '  +0.01%  '
$ws.Range("E36").Value = '  +1.98%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '164.06'
$ws.Range("E37").Value = '  +0.78%  '
$ws.Range("B38").Value = 'Mantle'
$ws.Range("C38").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.884'
$ws.Range("E38").Value = '  -1.44%  '
$ws.Range("B39").Value = 'Stacks'
$ws.Range("C39").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.95'
$ws.Range("E39").Value = '  +2.05%  '
$ws.Range("E40").Value = '  +3.71%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.83'
$ws.Range("E41").Value = '  +7.12%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '27.37'
$ws.Range("E42").Value = '  +3.42%  '
$ws.Range("E43").Value = '  +1.62%  '
$ws.Range("D44").Value = '2.911.80'
$ws.Range("E44").Value = '  +3.81%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '27.66'
$ws.Range("E45").Value = '  -0.25%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0744'
$ws.Range("E46").Value = '  -0.13%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '42.61'
$ws.Range("E47").Value = '  -0.73%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '353.11'
$ws.Range("E48").Value = '  +4.82%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0307'
$ws.Range("E49").Value = '  +0.48%  '
$ws.Range("E50").Value = '  -0.32%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '33.80'
$ws.Range("E51").Value = '  +0.45%  '
